$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.842016458511353
$ws.Range("B1").Value = 4.911418914794922
$ws.Range("C1").Value = 3.730206727981567
$ws.Range("D1").Value = 1.243616342544556
$ws.Range("E1").Value = 0.8185868263244629
